# Delete the "TOTAL" column (column C) on Sheet1 - it held combined
# "State, XX" strings that are no longer needed now that columns A (State)
# and B (State HEADLINES / abbreviation) already carry that information.
# This shifts the old column D ("Effective Tax Rate") left to become the
# new column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C:C").Delete()

# Move the active selection to match the workbook's saved cursor position.
[void]$ws.Range("F5").Select()
